$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 127, shifting existing rows 127:132 down to 128:133
$ws.Rows.Item(127).Insert()

# Fill the new row 127 with data (same as surrounding Puerro entries, but new date/price)
$ws.Cells.Item(127, 1).Value = 9
$ws.Cells.Item(127, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(127, 3).Value = "Metropolitana"
$ws.Cells.Item(127, 4).Value = 45147
$ws.Cells.Item(127, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(127, 5).Value = 13
$ws.Cells.Item(127, 6).Value = 100112005
$ws.Cells.Item(127, 7).Value = "Puerro"
$ws.Cells.Item(127, 8).Value = "Sin especificar"
$ws.Cells.Item(127, 9).Value = "Primera"
$ws.Cells.Item(127, 10).Value = 70
$ws.Cells.Item(127, 11).Value = 7000
$ws.Cells.Item(127, 12).Value = 7000
$ws.Cells.Item(127, 13).Value = 7000
$ws.Cells.Item(127, 14).Value = "`$/paquete 20 unidades"
$ws.Cells.Item(127, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(127, 16).Value = 350
$ws.Cells.Item(127, 17).Value = 20
$ws.Cells.Item(127, 18).Value = "Hortaliza"
